{"js": "// Load all paragraphs in the document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// --- Paragraph 1: \"Hello Tulane Technical Writing Students!\" -----------\n// The original paragraph is split into two runs bracketed by gramStart/\n// gramEnd proofing-error markers (\"Hello Tulane\" | \" Technical Writing\n// Students!\"). Re-writing the paragraph's text collapses it back into a\n// single run and drops the (now stale) grammar-check bookmarks.\nconst helloPara = paragraphs.items[0];\nhelloPara.clear();\nawait context.sync();\nhelloPara.insertText(\"Hello Tulane Technical Writing Students!\", Word.InsertLocation.start);\nawait context.sync();\n\n// --- Paragraph 5: \"To make sure you have ... write your name below:\" ---\n// Same situation: three runs split around the \"actually successfully\"\n// proofing-error span. Collapse to one run with the full sentence.\nconst confirmPara = paragraphs.items[4];\nconfirmPara.clear();\nawait context.sync();\nconfirmPara.insertText(\n  \"To make sure you have actually successfully viewed, edited, and pushed this document, please write your name below:\",\n  Word.InsertLocation.start\n);\nawait context.sync();\n\n// --- Paragraph 6: the bookmark (\"_GoBack\") paragraph --------------------\n// Insert the student's name as a new run *before* the existing\n// bookmarkStart/bookmarkEnd pair that already lives in that paragraph.\nconst namePara = paragraphs.items[5];\nnamePara.insertText(\"Kristen J. Rials (11/12/19)\", Word.InsertLocation.start);\nawait context.sync();\n\n// --- Trailing empty paragraph -------------------------------------------\n// The document gains one additional empty paragraph right after the\n// bookmark paragraph (and before the section properties).\nnamePara.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Word COM-interop (PowerShell-style) edit matching the same diff that\n# edit.js applies, expressed with $word.ActiveDocument / $d.Paragraphs /\n# Range methods.\n\n$d = $word.ActiveDocument\n\n# --- Paragraph 1: \"Hello Tulane Technical Writing Students!\" -----------\n# Originally split into two runs bracketed by gramStart/gramEnd proofing\n# marks (\"Hello Tulane\" | \" Technical Writing Students!\"). Deleting the\n# whole paragraph range (including its paragraph mark) removes the runs\n# *and* the now-orphaned w:proofErr markers in one shot; re-inserting a\n# fresh paragraph in its place with a single run gives the clean merged\n# text the diff expects.\n$p1 = $d.Paragraphs.Item(1)\n$p1.Range.Delete()\n$d.Paragraphs.Item(1).Range.InsertParagraphBefore()\n$d.Paragraphs.Item(1).Range.InsertBefore(\"Hello Tulane Technical Writing Students!\")\n\n# --- Paragraph 5: \"To make sure you have ... write your name below:\" ---\n# Same situation around the \"actually successfully\" proofing-error span;\n# collapse the three runs back into a single clean run.\n$p5 = $d.Paragraphs.Item(5)\n$p5.Range.Delete()\n$d.Paragraphs.Item(5).Range.InsertParagraphBefore()\n$d.Paragraphs.Item(5).Range.InsertBefore(\"To make sure you have actually successfully viewed, edited, and pushed this document, please write your name below:\")\n\n# --- Paragraph 6: the bookmark (\"_GoBack\") paragraph --------------------\n# Insert the student's name as a new run immediately before the existing\n# bookmarkStart/bookmarkEnd pair already in that paragraph.\n$p6 = $d.Paragraphs.Item(6)\n$p6.Range.InsertBefore(\"Kristen J. Rials (11/12/19)\")\n\n# --- Trailing empty paragraph -------------------------------------------\n# The document gains one additional empty paragraph right after the\n# bookmark paragraph (and before the section properties).\n$d.Paragraphs.Item(6).Range.InsertParagraphAfter()\n"}
